$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain plain text so numeric-looking
# strings (e.g. "1.001", "0.4884") are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.808.03"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "1.938.78"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "242.82"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "0.4884"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "0.2928"
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("D9").Value = "0.06900"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").Value = "19.17"
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("D11").Value = "104.78"
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D12").Value = "1.949.40"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "0.07715"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "5.351"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "0.6969"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "272.08"
$ws.Range("E16").Value = "  -5.04%  "
$ws.Range("D17").Value = "30.817.39"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "0.000007694"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "13.05"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.185.32"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "5.519"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "6.529"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("D25").Value = "9.700"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").Value = "166.24"
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("D27").Value = "19.53"
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("D28").Value = "2.155"
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("D29").Value = "0.1035"
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("E30").Value = "  -3.49%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.553"
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.552"
$ws.Range("E32").Value = "  -3.75%  "
$ws.Range("D33").Value = "4.352"
$ws.Range("E33").Value = "  -2.40%  "
$ws.Range("D34").Value = "0.04851"
$ws.Range("E34").Value = "  -3.16%  "
$ws.Range("D35").Value = "0.7532"
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("D36").Value = "1.152"
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("D37").Value = "1.000"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").Value = "0.01993"
$ws.Range("E39").Value = "  -2.32%  "
$ws.Range("D40").Value = "2.661"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("D41").Value = "6.476"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").Value = "77.11"
$ws.Range("E42").Value = "  +6.16%  "
$ws.Range("D43").Value = "2.078"
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("D44").Value = "0.9041"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("D45").Value = "0.4392"
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("D46").Value = "107.60"
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("D47").Value = "0.9991"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").Value = "7.677"
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("D49").Value = "982.60"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").Value = "0.1243"
$ws.Range("E50").Value = "  -3.00%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "9.273"
$ws.Range("E51").Value = "  -0.84%  "
